$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-5
# from serial 45174 (2023-09-05) to 45175 (2023-09-06)
$ws.Range("C2:C5").Value = 45175
